# Add "test 1" / "test 2" / "test 3" columns (E, F, G) as additional Jasmine
# test cases alongside the existing "initial values" column (D), per the
# commit "add additional calcs to provide test cases to jasmine".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rightAlign  = -4152
$centerAlign = -4108

$currencyCentsFmt = '_("$"* #,##0.00_);_("$"* \(#,##0.00\);_("$"* "-"??_);_(@_)'
$currencyWholeFmt = '_("$"* #,##0_);_("$"* \(#,##0\);_("$"* "-"??_);_(@_)'

# --- Row 2: header labels (D2 "initial values" reuses column D's default
#     right-aligned style; E2:G2 reuse the centered style used elsewhere in
#     the sheet) -----------------------------------------------------------
$ws.Range("D2").Value = "initial values"
$ws.Range("D2").HorizontalAlignment = $rightAlign

$ws.Range("E2").Value = "test 1"
$ws.Range("E2").HorizontalAlignment = $centerAlign

$ws.Range("F2").Value = "test 2"
$ws.Range("F2").HorizontalAlignment = $centerAlign

$ws.Range("G2").Value = "test 3"
$ws.Range("G2").HorizontalAlignment = $centerAlign

# --- Row 3: Result formulas (currency, 2dp, right aligned) ----------------
$ws.Range("E3").Formula = "=(E5*E6)/(1-(1+E6)^(-E7))"
$ws.Range("E3").HorizontalAlignment = $rightAlign
$ws.Range("E3").NumberFormat = $currencyCentsFmt

$ws.Range("F3").Formula = "=(F5*F6)/(1-(1+F6)^(-F7))"
$ws.Range("F3").HorizontalAlignment = $rightAlign
$ws.Range("F3").NumberFormat = $currencyCentsFmt

$ws.Range("G3").Formula = "=(G5*G6)/(1-(1+G6)^(-G7))"
$ws.Range("G3").HorizontalAlignment = $rightAlign
$ws.Range("G3").NumberFormat = $currencyCentsFmt

# --- Row 4: spacer row, right-aligned blank cells --------------------------
$ws.Range("E4:G4").HorizontalAlignment = $rightAlign

# --- Row 5: Principle formulas (=row10, currency whole $, right) ----------
$ws.Range("E5").Formula = "=E10"
$ws.Range("E5").HorizontalAlignment = $rightAlign
$ws.Range("E5").NumberFormat = $currencyWholeFmt

$ws.Range("F5").Formula = "=F10"
$ws.Range("F5").HorizontalAlignment = $rightAlign
$ws.Range("F5").NumberFormat = $currencyWholeFmt

$ws.Range("G5").Formula = "=G10"
$ws.Range("G5").HorizontalAlignment = $rightAlign
$ws.Range("G5").NumberFormat = $currencyWholeFmt

# --- Row 6: periodic interest rate formulas (centered) ---------------------
$ws.Range("E6").Formula = "=E12/100/12"
$ws.Range("E6").HorizontalAlignment = $centerAlign

$ws.Range("F6").Formula = "=F12/100/12"
$ws.Range("F6").HorizontalAlignment = $centerAlign

$ws.Range("G6").Formula = "=G12/100/12"
$ws.Range("G6").HorizontalAlignment = $centerAlign

# --- Row 7: total # of payments formulas (centered) ------------------------
$ws.Range("E7").Formula = "=E11*12"
$ws.Range("E7").HorizontalAlignment = $centerAlign

$ws.Range("F7").Formula = "=F11*12"
$ws.Range("F7").HorizontalAlignment = $centerAlign

$ws.Range("G7").Formula = "=G11*12"
$ws.Range("G7").HorizontalAlignment = $centerAlign

# --- Row 8 & 9: spacer rows, right-aligned blank cells ----------------------
$ws.Range("E8:G8").HorizontalAlignment = $rightAlign
$ws.Range("E9:G9").HorizontalAlignment = $rightAlign

# --- Row 10: Loan amount inputs (currency, right) ---------------------------
$ws.Range("E10").Value = 100000
$ws.Range("E10").HorizontalAlignment = $rightAlign
$ws.Range("E10").NumberFormat = $currencyWholeFmt

$ws.Range("F10").Value = 500000
$ws.Range("F10").HorizontalAlignment = $rightAlign
$ws.Range("F10").NumberFormat = $currencyWholeFmt

$ws.Range("G10").Value = 12500.5
$ws.Range("G10").HorizontalAlignment = $rightAlign
$ws.Range("G10").NumberFormat = $currencyCentsFmt

# --- Row 11: Term in Years inputs (centered) --------------------------------
$ws.Range("E11").Value = 15
$ws.Range("E11").HorizontalAlignment = $centerAlign

$ws.Range("F11").Value = 10
$ws.Range("F11").HorizontalAlignment = $centerAlign

$ws.Range("G11").Value = 20
$ws.Range("G11").HorizontalAlignment = $centerAlign

# --- Row 12: Yearly Rate inputs (centered) ----------------------------------
$ws.Range("E12").Value = 3
$ws.Range("E12").HorizontalAlignment = $centerAlign

$ws.Range("F12").Value = 4.785
$ws.Range("F12").HorizontalAlignment = $centerAlign

$ws.Range("G12").Value = 0.0001
$ws.Range("G12").HorizontalAlignment = $centerAlign

# --- Column widths for the new test columns ---------------------------------
$ws.Columns("E:F").ColumnWidth = 10
$ws.Columns("G").ColumnWidth = 11.5

# --- Selection moves to G13, matching the saved workbook state -------------
$ws.Range("G13").Select()

Write-Output "edit applied"
